$d = $word.ActiveDocument

# --- Table 1: Grundpreis / Arbeitspreis updates ---
$t = $d.Tables.Item(1)

# Row 2 ("Watt für wenig"): Grundgebühr 8,20 € -> 13,50 €, Verbrauchspreis 0,16 € -> 0,75 €
$t.Cell(2,2).Range.Find.Execute("8,20 €", $true, $false, $false, $false, $false, $true, 1, $false, "13,50 €", 2)
$t.Cell(2,3).Range.Find.Execute("0,16 €", $true, $false, $false, $false, $false, $true, 1, $false, "0,75 €", 2)

# Row 3 ("Billig Strom"): Grundgebühr 4,90 € -> 9,20 €, Verbrauchspreis 0,19 € -> 0,81 €
$t.Cell(3,2).Range.Find.Execute("4,90 €", $true, $false, $false, $false, $false, $true, 1, $false, "9,20 €", 2)
$t.Cell(3,3).Range.Find.Execute("0,19 €", $true, $false, $false, $false, $false, $true, 1, $false, "0,81 €", 2)

# --- Header: SAVEDATE field cached result ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$rng = $hdr.Range.Duplicate
$rng.Find.Execute("04.09.2022", $true, $false, $false, $false, $false, $true, 1, $false, "06.10.2022", 2)
